$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Column H: Chinese move name (chName) ---
# H1 is the header cell; copy the bold header style used by A1:G1
# so H1 matches the look of the rest of the header row.
$ws.Range("A1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "chName"

# H2:H121 hold the localized (Chinese) move names, one per data row.
# (Row 111 intentionally has no chName value, matching the source data.)
$ws.Range("H2").Value = "超级角击"
$ws.Range("H3").Value = "十字剪"
$ws.Range("H4").Value = "信号光束"
$ws.Range("H5").Value = "银色旋风"
$ws.Range("H6").Value = "虫鸣"
$ws.Range("H7").Value = "欺诈"
$ws.Range("H8").Value = "恶之波动"
$ws.Range("H9").Value = "暗袭要害"
$ws.Range("H10").Value = "咬碎"
$ws.Range("H11").Value = "龙之波动"
$ws.Range("H12").Value = "逆鳞"
$ws.Range("H13").Value = "龙爪"
$ws.Range("H14").Value = "龙卷风"
$ws.Range("H15").Value = "打雷"
$ws.Range("H16").Value = "抛物面充电"
$ws.Range("H17").Value = "雷电拳"
$ws.Range("H18").Value = "疯狂伏特"
$ws.Range("H19").Value = "放电"
$ws.Range("H20").Value = "电磁炮"
$ws.Range("H21").Value = "十万伏特"
$ws.Range("H22").Value = "月亮之力"
$ws.Range("H23").Value = "魔法闪耀"
$ws.Range("H24").Value = "吸取之吻"
$ws.Range("H25").Value = "嬉闹"
$ws.Range("H26").Value = "魅惑之声"
$ws.Range("H27").Value = "近身战"
$ws.Range("H28").Value = "地狱翻滚"
$ws.Range("H29").Value = "十字劈"
$ws.Range("H30").Value = "下盘踢"
$ws.Range("H31").Value = "真气弹"
$ws.Range("H32").Value = "劈瓦"
$ws.Range("H33").Value = "爆裂拳"
$ws.Range("H34").Value = "蓄能焰袭"
$ws.Range("H35").Value = "火焰拳"
$ws.Range("H36").Value = "大字爆炎"
$ws.Range("H37").Value = "过热"
$ws.Range("H38").Value = "火焰轮"
$ws.Range("H39").Value = "烈焰溅射"
$ws.Range("H40").Value = "喷射火焰"
$ws.Range("H41").Value = "热风"
$ws.Range("H42").Value = "勇鸟猛攻"
$ws.Range("H43").Value = "神鸟猛攻"
$ws.Range("H44").Value = "燕返"
$ws.Range("H45").Value = "暴风"
$ws.Range("H46").Value = "空气利刃"
$ws.Range("H47").Value = "啄钻"
$ws.Range("H48").Value = "黑夜魔影"
$ws.Range("H49").Value = "暗影拳"
$ws.Range("H50").Value = "影子偷袭"
$ws.Range("H51").Value = "暗影球"
$ws.Range("H52").Value = "奇异之风"
$ws.Range("H53").Value = "种子炸弹"
$ws.Range("H54").Value = "超级吸取"
$ws.Range("H55").Value = "日光束"
$ws.Range("H56").Value = "打草结"
$ws.Range("H57").Value = "蛮力藤鞭"
$ws.Range("H58").Value = "落英缤纷"
$ws.Range("H59").Value = "能量球"
$ws.Range("H60").Value = "叶刃"
$ws.Range("H61").Value = "终极吸取"
$ws.Range("H62").Value = "地震"
$ws.Range("H63").Value = "重踏"
$ws.Range("H64").Value = "骨棒"
$ws.Range("H65").Value = "泥巴炸弹"
$ws.Range("H66").Value = "流沙地狱"
$ws.Range("H67").Value = "直冲钻"
$ws.Range("H68").Value = "挖洞"
$ws.Range("H69").Value = "极光束"
$ws.Range("H70").Value = "冰冻光线"
$ws.Range("H71").Value = "冰冻拳"
$ws.Range("H72").Value = "雪崩"
$ws.Range("H73").Value = "冰冻之风"
$ws.Range("H74").Value = "暴风雪"
$ws.Range("H75").Value = "夹住"
$ws.Range("H76").Value = "破坏光线"
$ws.Range("H77").Value = "紧束"
$ws.Range("H78").Value = "必杀门牙"
$ws.Range("H79").Value = "高速星星"
$ws.Range("H80").Value = "挣扎"
$ws.Range("H81").Value = "角撞"
$ws.Range("H82").Value = "睡觉"
$ws.Range("H83").Value = "泰山压顶"
$ws.Range("H84").Value = "踩踏"
$ws.Range("H85").Value = "剧毒牙"
$ws.Range("H86").Value = "污泥炸弹"
$ws.Range("H87").Value = "污泥波"
$ws.Range("H88").Value = "污泥攻击"
$ws.Range("H89").Value = "垃圾射击"
$ws.Range("H90").Value = "十字毒刃"
$ws.Range("H91").Value = "爱心印章"
$ws.Range("H92").Value = "精神击破"
$ws.Range("H93").Value = "预知未来"
$ws.Range("H94").Value = "镜面反射"
$ws.Range("H95").Value = "精神冲击"
$ws.Range("H96").Value = "幻象光线"
$ws.Range("H97").Value = "精神强念"
$ws.Range("H98").Value = "岩崩"
$ws.Range("H99").Value = "岩石爆击"
$ws.Range("H100").Value = "原始之力"
$ws.Range("H101").Value = "尖石攻击"
$ws.Range("H102").Value = "岩石封锁"
$ws.Range("H103").Value = "力量宝石"
$ws.Range("H104").Value = "磁铁炸弹"
$ws.Range("H105").Value = "铁头"
$ws.Range("H106").Value = "重磅冲撞"
$ws.Range("H107").Value = "陀螺球"
$ws.Range("H108").Value = "加农光炮"
$ws.Range("H109").Value = "泡沫光线"
$ws.Range("H110").Value = "水炮"
$ws.Range("H112").Value = "水流喷射"
$ws.Range("H113").Value = "盐水"
$ws.Range("H114").Value = "水之波动"
$ws.Range("H115").Value = "热水"
$ws.Range("H116").Value = "断崖之剑"
$ws.Range("H117").Value = "破灭之愿"
$ws.Range("H118").Value = "流星群"
$ws.Range("H119").Value = "精神突进"
$ws.Range("H120").Value = "根源波动"
$ws.Range("H121").Value = "冲浪"

# Restore the workbook selection/cursor position to match the edited view
$ws.Range("H105").Select() | Out-Null
